$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and two row re-orderings)

# Row 2
$ws.Range("D2").Value = '66.148.20'
$ws.Range("E2").Value = '  -1.57%  '

# Row 3
$ws.Range("D3").Value = '3.438.62'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.39'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.55%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.91'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.35%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.09%  '

# Row 9
$ws.Range("D9").Value = '3.433.43'
$ws.Range("E9").Value = '  -1.12%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.131'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.61%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.85'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.11%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.418'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.38%  '

# Row 13
$ws.Range("D13").Value = '4.041.06'
$ws.Range("E13").Value = '  -0.93%  '

# Row 14
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.132'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.47%  '

# Row 15
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.97'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -5.44%  '

# Row 16
$ws.Range("D16").Value = '66.209.52'
$ws.Range("E16").Value = '  -1.60%  '

# Row 17
$ws.Range("E17").Value = '  -3.29%  '

# Row 18
$ws.Range("D18").Value = '3.453.01'
$ws.Range("E18").Value = '  -1.00%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.97'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.92%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.68'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.34%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '377.32'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.80%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.73'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.23%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.542'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.90%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.03%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '71.92'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.07%  '

# Row 26
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.72'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.66%  '

# Row 27
$ws.Range("E27").Value = '  -3.72%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.76'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.91%  '

# Row 29
$ws.Range("E29").Value = '  -1.06%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.17%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '24.08'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.70%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.81'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -5.34%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.98'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.84%  '

# Row 34
$ws.Range("E34").Value = '  +0.02%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.30'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -7.69%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.10'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.49%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.56'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.67%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '159.21'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.68%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '28.98'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +10.92%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.881'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.00%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.76'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -5.58%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.47'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.80%  '

# Row 43
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.52'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -10.86%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.35'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.74%  '

# Row 45
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0690'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.32%  '

# Row 46
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.657.73'
$ws.Range("E46").Value = '  -5.28%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.26'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.31%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.20'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -9.01%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0288'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.69%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '310.75'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -6.18%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.998'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -5.05%  '
